$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 423.78946
$ws.Range("I33").Value = 272.06668
$ws.Range("J33").Value = 992.75
$ws.Range("K33").Value = 272.06668
$ws.Range("L33").Value = 992.75
$ws.Range("M33").Value = -43.06668000000002
$ws.Range("N33").Value = -1450.75
$ws.Range("H40").Value = 94486.17999999999
$ws.Range("J40").Value = 3881.8928
$ws.Range("L40").Value = 3881.8928
$ws.Range("N40").Value = -4231.8928
$ws.Range("H98").Value = 4775.9287
$ws.Range("I98").Value = 5105.6665
$ws.Range("K98").Value = 5105.6665
$ws.Range("M98").Value = -3607.6665
$ws.Range("H122").Value = 4775.9287
$ws.Range("I122").Value = 5105.6665
$ws.Range("K122").Value = 15316.9995
$ws.Range("M122").Value = -12866.9995
$ws.Range("H128").Value = 99894
$ws.Range("J128").Value = 99894
$ws.Range("L128").Value = 99894
$ws.Range("N128").Value = -109854
$ws.Range("H137").Value = 1413.4117
$ws.Range("I137").Value = 1416.0834
$ws.Range("J137").Value = 1407
$ws.Range("K137").Value = 4248.2502
$ws.Range("L137").Value = 4221
$ws.Range("M137").Value = -1698.2502
$ws.Range("N137").Value = -9321
$ws.Range("H138").Value = 4266.41
$ws.Range("J138").Value = 4426.573
$ws.Range("L138").Value = 13279.719
$ws.Range("N138").Value = -23559.719

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31717936
$ws.Range("I32").Value = 31747122
$ws.Range("K32").Value = 31747122
$ws.Range("M32").Value = -31746835
$ws.Range("H74").Value = 1295.5714
$ws.Range("I74").Value = 1231.25
$ws.Range("K74").Value = 1231.25
$ws.Range("M74").Value = -357.25
$ws.Range("H77").Value = 1295.5714
$ws.Range("I77").Value = 1231.25
$ws.Range("K77").Value = 6156.25
$ws.Range("M77").Value = -1788.25
$ws.Range("H97").Value = 1793.1666
$ws.Range("I97").Value = 1269.7858
$ws.Range("K97").Value = 1269.7858
$ws.Range("M97").Value = -773.7858000000001
$ws.Range("H122").Value = 2669.423
$ws.Range("I122").Value = 2391.9167
$ws.Range("K122").Value = 7175.750100000001
$ws.Range("M122").Value = -4725.750100000001
$ws.Range("H132").Value = 458609.03
$ws.Range("I132").Value = 717374.9
$ws.Range("K132").Value = 2152124.7
$ws.Range("M132").Value = -2149594.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 296.75
$ws.Range("I22").Value = 296.75
$ws.Range("K22").Value = 296.75
$ws.Range("M22").Value = -123.75
$ws.Range("H105").Value = 3309.8
$ws.Range("I105").Value = 3699.7144
$ws.Range("J105").Value = 2400
$ws.Range("K105").Value = 3699.7144
$ws.Range("L105").Value = 2400
$ws.Range("M105").Value = -1952.7144
$ws.Range("N105").Value = -5894
$ws.Range("H134").Value = 1964034.6
$ws.Range("I134").Value = 2383449
$ws.Range("K134").Value = 7150347
$ws.Range("M134").Value = -7147812

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2497.0378
$ws.Range("I31").Value = 1459.6
$ws.Range("J31").Value = 3423.3215
$ws.Range("K31").Value = 1459.6
$ws.Range("L31").Value = 3423.3215
$ws.Range("M31").Value = -1164.6
$ws.Range("N31").Value = -4013.3215
$ws.Range("H34").Value = 2497.0378
$ws.Range("I34").Value = 1459.6
$ws.Range("J34").Value = 3423.3215
$ws.Range("K34").Value = 1459.6
$ws.Range("L34").Value = 3423.3215
$ws.Range("M34").Value = -1257.6
$ws.Range("N34").Value = -3827.3215
$ws.Range("H58").Value = 3862.1143
$ws.Range("I58").Value = 2963.2307
$ws.Range("J58").Value = 4393.273
$ws.Range("K58").Value = 2963.2307
$ws.Range("L58").Value = 4393.273
$ws.Range("M58").Value = -2760.2307
$ws.Range("N58").Value = -4799.273
$ws.Range("H105").Value = 1989.9445
$ws.Range("I105").Value = 1921.4667
$ws.Range("J105").Value = 2332.3333
$ws.Range("K105").Value = 1921.4667
$ws.Range("L105").Value = 2332.3333
$ws.Range("M105").Value = -174.4666999999999
$ws.Range("N105").Value = -5826.3333
$ws.Range("H132").Value = 3363.5312
$ws.Range("I132").Value = 3252.0715
$ws.Range("K132").Value = 9756.2145
$ws.Range("M132").Value = -7226.2145
$ws.Range("H136").Value = 3862.1143
$ws.Range("I136").Value = 2963.2307
$ws.Range("J136").Value = 4393.273
$ws.Range("K136").Value = 8889.6921
$ws.Range("L136").Value = 13179.819
$ws.Range("M136").Value = -6339.6921
$ws.Range("N136").Value = -18279.819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 50500444
$ws.Range("J11").Value = 83333864
$ws.Range("L11").Value = 250001592
$ws.Range("N11").Value = -250001872
$ws.Range("H13").Value = 18.333334
$ws.Range("J13").Value = 30
$ws.Range("L13").Value = 90
$ws.Range("N13").Value = -426
$ws.Range("H16").Value = 199.66667
$ws.Range("I16").Value = 199.66667
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 599.00001
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -426.00001
$ws.Range("N16").ClearContents()
$ws.Range("H23").Value = 522.8333
$ws.Range("I23").Value = 570.6667
$ws.Range("J23").Value = 475
$ws.Range("K23").Value = 1712.0001
$ws.Range("L23").Value = 1425
$ws.Range("M23").Value = -1477.0001
$ws.Range("N23").Value = -1895
$ws.Range("H29").Value = 2873
$ws.Range("I29").Value = 3732.6667
$ws.Range("J29").Value = 294
$ws.Range("K29").Value = 11198.0001
$ws.Range("L29").Value = 882
$ws.Range("M29").Value = -10921.0001
$ws.Range("N29").Value = -1436
$ws.Range("H131").Value = 2002.6072
$ws.Range("I131").Value = 2643.5557
$ws.Range("J131").Value = 1699
$ws.Range("K131").Value = 7930.6671
$ws.Range("L131").Value = 5097
$ws.Range("M131").Value = -2890.6671
$ws.Range("N131").Value = -15177
$ws.Range("H139").Value = 1793.5714
$ws.Range("I139").Value = 1793.5714
$ws.Range("K139").Value = 5380.7142
$ws.Range("M139").Value = -240.7142000000003
$ws.Range("H140").Value = 1086.7273
$ws.Range("I140").Value = 1006
$ws.Range("K140").Value = 3018
$ws.Range("M140").Value = 2162

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1894.878
$ws.Range("I97").Value = 1213.0322
$ws.Range("K97").Value = 1213.0322
$ws.Range("M97").Value = -717.0322000000001
$ws.Range("H102").Value = 2553.1155
$ws.Range("I102").Value = 2392.8635
$ws.Range("J102").Value = 3434.5
$ws.Range("K102").Value = 2392.8635
$ws.Range("L102").Value = 3434.5
$ws.Range("M102").Value = -770.8634999999999
$ws.Range("N102").Value = -6678.5
$ws.Range("H113").Value = 20761.375
$ws.Range("I113").Value = 3432.7222
$ws.Range("J113").Value = 34939.363
$ws.Range("K113").Value = 3432.7222
$ws.Range("L113").Value = 34939.363
$ws.Range("M113").Value = -1262.7222
$ws.Range("N113").Value = -39279.363
$ws.Range("H122").Value = 1386.3684
$ws.Range("I122").Value = 1515.0625
$ws.Range("J122").Value = 700
$ws.Range("K122").Value = 4545.1875
$ws.Range("L122").Value = 2100
$ws.Range("M122").Value = -2095.1875
$ws.Range("N122").Value = -7000
$ws.Range("H135").Value = 83878
$ws.Range("J135").Value = 83878
$ws.Range("L135").Value = 83878
$ws.Range("N135").Value = -94018

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9249.75
$ws.Range("I7").Value = 8999.666999999999
$ws.Range("K7").Value = 8999.666999999999
$ws.Range("M7").Value = -8887.666999999999
$ws.Range("H40").Value = 3070.6453
$ws.Range("I40").Value = 3023.1482
$ws.Range("K40").Value = 3023.1482
$ws.Range("M40").Value = -2887.1482
$ws.Range("H114").Value = 95211.25
$ws.Range("J114").Value = 95211.25
$ws.Range("L114").Value = 95211.25
$ws.Range("N114").Value = -103889.25
$ws.Range("H122").Value = 10524.611
$ws.Range("I122").Value = 9498.308000000001
$ws.Range("J122").Value = 13193
$ws.Range("K122").Value = 28494.924
$ws.Range("L122").Value = 39579
$ws.Range("M122").Value = -26044.924
$ws.Range("N122").Value = -44479
$ws.Range("H126").Value = 9249.75
$ws.Range("I126").Value = 8999.666999999999
$ws.Range("K126").Value = 26999.001
$ws.Range("M126").Value = -24529.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 69500
$ws.Range("J16").Value = 69500
$ws.Range("L16").Value = 69500
$ws.Range("N16").Value = -70084
$ws.Range("H62").Value = 8398.875
$ws.Range("J62").Value = 6140
$ws.Range("L62").Value = 6140
$ws.Range("N62").Value = -7388
$ws.Range("H65").Value = 8398.875
$ws.Range("J65").Value = 6140
$ws.Range("L65").Value = 30700
$ws.Range("N65").Value = -36940
$ws.Range("H81").Value = 72937.53
$ws.Range("J81").Value = 9999.833000000001
$ws.Range("L81").Value = 19999.666
$ws.Range("N81").Value = -22121.666
$ws.Range("H84").Value = 72937.53
$ws.Range("J84").Value = 9999.833000000001
$ws.Range("L84").Value = 99998.33
$ws.Range("N84").Value = -110606.33
$ws.Range("H122").Value = 5710.7915
$ws.Range("I122").Value = 5436.9443
$ws.Range("J122").Value = 6532.3335
$ws.Range("K122").Value = 16310.8329
$ws.Range("L122").Value = 19597.0005
$ws.Range("M122").Value = -13860.8329
$ws.Range("N122").Value = -24497.0005
$ws.Range("H132").Value = 19254.707
$ws.Range("I132").Value = 22532.857
$ws.Range("K132").Value = 67598.571
$ws.Range("M132").Value = -65068.571
